# last changes to v1.8.2
$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Sheets.Item("Metadata")
$wsMeta.Range("B3").Value = "1.8.2"
$wsMeta.Range("B8").Value = "2023-09-01T14:45:29-04:00"

$wsElements = $wb.Sheets.Item("Elements")
$wsElements.Range("AJ1").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}
ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
